# UndoRedoActivityDiagram.pptx edit:
#  - "The Food Diary" re-theme: "address book" -> "food diary" wording
#    (and the matching "addressBookStateList" -> "foodDiaryStateList"
#    variable-name text) on slide 1.
#  - Refreshed "datetimeFigureOut" footer field cache text (6/7/2018 ->
#    5/4/19) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

function Update-DateShape($shape) {
    if (-not $shape.HasTextFrame) { return }
    $tr = $shape.TextFrame.TextRange
    if ($tr.Text -eq "6/7/2018") {
        $tr.Text = "5/4/19"
    }
}

# --- Slide master: refresh the cached date field text ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

# --- Every slide layout: refresh the cached date field text ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# --- Slide 1: rewrite the two "address book" text boxes to "food diary" ---
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if (-not $shape.HasTextFrame) { continue }
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text

    if ($full -eq "[command commits address book]") {
        # Keep the leading "[" as its own run; rewrite the rest of the
        # run ("command commits address book]") in one go so it stays a
        # single run, matching the original run split.
        $oldRun = "command commits address book]"
        $newRun = "command commits food diary]"
        $idx = $full.IndexOf($oldRun)
        $c = $tr.Characters($idx + 1, $oldRun.Length)
        $c.Text = $newRun
    }
    elseif ($full -eq "Purge redundant states and then save address book to addressBookStateList ") {
        $oldRun0 = "Purge redundant states and then save address book to "
        $newRun0 = "Purge redundant states and then save food diary to "
        $idx0 = $full.IndexOf($oldRun0)
        $c0 = $tr.Characters($idx0 + 1, $oldRun0.Length)
        $c0.Text = $newRun0

        $full2 = $tr.Text
        $oldRun1 = "addressBookStateList"
        $newRun1 = "foodDiaryStateList"
        $idx1 = $full2.IndexOf($oldRun1)
        $c1 = $tr.Characters($idx1 + 1, $oldRun1.Length)
        $c1.Text = $newRun1
    }
}
